$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a new paragraph "Play Buffalo Blitz Free: Review of
#    Unique 6x4 Online Slot Game" (bold) right before the final
#    "Prompt: ..." paragraph.
#
#    InsertXML only materialises a paragraph break *between* <w:p>
#    elements in the fragment, so we append a throw-away empty <w:p>
#    after our real paragraph to force the break, then delete the
#    throw-away paragraph that is left behind.
# ------------------------------------------------------------------
$paraCountBeforeInsert = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($paraCountBeforeInsert)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Buffalo Blitz Free: Review of Unique 6x4 Online Slot Game</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"></w:p>'
$insertPoint.InsertXML($newParaXml)

# Remove the left-over empty paragraph created only to force the break.
$spacerParaIndex = $paraCountBeforeInsert + 1
$d.Paragraphs($spacerParaIndex).Range.Delete()

# ------------------------------------------------------------------
# 2. Remove the old "Meta description" paragraph (the second
#    paragraph in the document, right after the title).
# ------------------------------------------------------------------
$d.Paragraphs(2).Range.Delete()

# ------------------------------------------------------------------
# 3. Replace the text of the (still-italic) "Prompt: ..." paragraph
#    with the former meta-description sentence.
# ------------------------------------------------------------------
$oldPrompt = 'Prompt: Design a feature image for "Buffalo Blitz" that showcases a happy Maya warrior with glasses in a cartoon style. The image should feature the warrior standing in the midst of the North American prairie, with various animals such as moose, raccoons, pumas, and bears around him. The warrior should be wearing traditional Maya clothing, including a headdress and a necklace made of buffalo bones. The glasses should be modern and stylish, to contrast with the traditional clothing. The image should have a bright and vibrant color scheme, with the warrior smiling while holding a buffalo horn, as if ready to start playing the game. The image should highlight the unique gameplay of Buffalo Blitz with a text overlay that reads "More Symbols, More Fun: Play Buffalo Blitz Now!"'
$newPrompt = 'Find out why Buffalo Blitz offers a dynamic and engaging experience with a 6x4 grid and 4096 paylines. Play free and read our review now.'

$d.Content.Find.Execute($oldPrompt, $false, $false, $false, $false, $false, $true, 1, $false, $newPrompt, 2)
